{"js": "// Office.js (Word JavaScript API) edit script\n// Replaces the 100 arithmetic-equation strings inside the single table's\n// cells with their updated counterparts, in reading order (row-major),\n// exactly as in the target diff. Run properties / formatting of each\n// cell are left untouched because only the table's text `values` are\n// rewritten, not the run XML itself.\n\nconst replacements = [\n  [\"32+63=95\", \"28-0=28\"],\n  [\"83-58=25\", \"80-56=24\"],\n  [\"39-15=24\", \"49+33=82\"],\n  [\"74-44=30\", \"9+2=11\"],\n  [\"27-12=15\", \"4+28=32\"],\n  [\"78-58=20\", \"55+8=63\"],\n  [\"19-15=4\", \"95-63=32\"],\n  [\"70-18=52\", \"46+48=94\"],\n  [\"49-28=21\", \"60-59=1\"],\n  [\"83-61=22\", \"34+60=94\"],\n  [\"49-7=42\", \"96-1=95\"],\n  [\"56-16=40\", \"59-32=27\"],\n  [\"81+13=94\", \"85-71=14\"],\n  [\"0+34=34\", \"62+35=97\"],\n  [\"21+69=90\", \"18+68=86\"],\n  [\"14+73=87\", \"70-35=35\"],\n  [\"15+63=78\", \"87-52=35\"],\n  [\"15+74=89\", \"28+20=48\"],\n  [\"90+6=96\", \"74-22=52\"],\n  [\"59+14=73\", \"0+83=83\"],\n  [\"20+20=40\", \"66-10=56\"],\n  [\"45+47=92\", \"3+41=44\"],\n  [\"77-11=66\", \"65+2=67\"],\n  [\"45-35=10\", \"98-66=32\"],\n  [\"79-76=3\", \"83-72=11\"],\n  [\"76+9=85\", \"3+0=3\"],\n  [\"14+78=92\", \"97-12=85\"],\n  [\"8+17=25\", \"3-3=0\"],\n  [\"9+34=43\", \"6+2=8\"],\n  [\"37-5=32\", \"11-5=6\"],\n  [\"67-61=6\", \"39+9=48\"],\n  [\"0+97=97\", \"73-28=45\"],\n  [\"56-6=50\", \"75-4=71\"],\n  [\"67-16=51\", \"6+57=63\"],\n  [\"48-12=36\", \"46+3=49\"],\n  [\"1+98=99\", \"11+83=94\"],\n  [\"37-22=15\", \"33-31=2\"],\n  [\"83-20=63\", \"89-23=66\"],\n  [\"40-35=5\", \"36+41=77\"],\n  [\"31-20=11\", \"12+50=62\"],\n  [\"58-52=6\", \"6+14=20\"],\n  [\"48+48=96\", \"92-79=13\"],\n  [\"14+15=29\", \"66-19=47\"],\n  [\"49-49=0\", \"97-37=60\"],\n  [\"37+62=99\", \"27+20=47\"],\n  [\"51+48=99\", \"42+31=73\"],\n  [\"37-13=24\", \"77-24=53\"],\n  [\"73-69=4\", \"21+23=44\"],\n  [\"74-55=19\", \"75+13=88\"],\n  [\"21-9=12\", \"94-43=51\"],\n  [\"63-50=13\", \"46+44=90\"],\n  [\"42-18=24\", \"57+42=99\"],\n  [\"96-2=94\", \"58-58=0\"],\n  [\"56-39=17\", \"84-0=84\"],\n  [\"84-14=70\", \"73-0=73\"],\n  [\"26+73=99\", \"70+29=99\"],\n  [\"65-45=20\", \"48+31=79\"],\n  [\"25+24=49\", \"52-46=6\"],\n  [\"53+39=92\", \"90-31=59\"],\n  [\"75-34=41\", \"70-39=31\"],\n  [\"96-36=60\", \"75-52=23\"],\n  [\"48+6=54\", \"97-43=54\"],\n  [\"99-64=35\", \"25+0=25\"],\n  [\"32+31=63\", \"21+16=37\"],\n  [\"99-78=21\", \"90-68=22\"],\n  [\"52-18=34\", \"43-11=32\"],\n  [\"35-33=2\", \"98-23=75\"],\n  [\"80-61=19\", \"28+28=56\"],\n  [\"88-43=45\", \"22+11=33\"],\n  [\"31+68=99\", \"45-1=44\"],\n  [\"72-26=46\", \"1+62=63\"],\n  [\"87-31=56\", \"56-38=18\"],\n  [\"43+3=46\", \"92-90=2\"],\n  [\"35+31=66\", \"92-40=52\"],\n  [\"76+7=83\", \"11+68=79\"],\n  [\"50+32=82\", \"76-3=73\"],\n  [\"23+7=30\", \"43-27=16\"],\n  [\"74+6=80\", \"27-26=1\"],\n  [\"59-45=14\", \"11+42=53\"],\n  [\"17+17=34\", \"43+7=50\"],\n  [\"89-52=37\", \"69-20=49\"],\n  [\"5+70=75\", \"65+24=89\"],\n  [\"7+30=37\", \"64+25=89\"],\n  [\"27+47=74\", \"21-19=2\"],\n  [\"62+5=67\", \"71+0=71\"],\n  [\"65-46=19\", \"80-8=72\"],\n  [\"1+85=86\", \"51+26=77\"],\n  [\"86-6=80\", \"59-59=0\"],\n  [\"92-67=25\", \"71-19=52\"],\n  [\"13+32=45\", \"52+33=85\"],\n  [\"43-35=8\", \"52+19=71\"],\n  [\"80-3=77\", \"31+30=61\"],\n  [\"68-55=13\", \"37+23=60\"],\n  [\"82+2=84\", \"65-12=53\"],\n  [\"19+72=91\", \"80-45=35\"],\n  [\"78-15=63\", \"72-71=1\"],\n  [\"13+37=50\", \"38+51=89\"],\n  [\"89+10=99\", \"96-25=71\"],\n  [\"38+10=48\", \"93-56=37\"],\n  [\"67-13=54\", \"61+18=79\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values,rowCount,columnCount\");\nawait context.sync();\n\nconst values = table.values;\nlet idx = 0;\nfor (let r = 0; r < values.length; r++) {\n  for (let c = 0; c < values[r].length; c++) {\n    if (idx >= replacements.length) break;\n    const [oldText, newText] = replacements[idx];\n    if (values[r][c] === oldText) {\n      values[r][c] = newText;\n    }\n    idx++;\n  }\n}\n\ntable.values = values;\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell) edit script\n# Replaces the 100 arithmetic-equation strings inside the single table's\n# cells with their updated counterparts. Pairs are matched positionally,\n# in table reading order (row-major), against the content captured in\n# the target diff; each cell's current text is also verified against the\n# expected \"old\" value before being overwritten, so formatting (the\n# run's rPr) is left completely untouched -- only the text changes.\n\n$replacements = @(\n  @(\"32+63=95\", \"28-0=28\"),\n  @(\"83-58=25\", \"80-56=24\"),\n  @(\"39-15=24\", \"49+33=82\"),\n  @(\"74-44=30\", \"9+2=11\"),\n  @(\"27-12=15\", \"4+28=32\"),\n  @(\"78-58=20\", \"55+8=63\"),\n  @(\"19-15=4\", \"95-63=32\"),\n  @(\"70-18=52\", \"46+48=94\"),\n  @(\"49-28=21\", \"60-59=1\"),\n  @(\"83-61=22\", \"34+60=94\"),\n  @(\"49-7=42\", \"96-1=95\"),\n  @(\"56-16=40\", \"59-32=27\"),\n  @(\"81+13=94\", \"85-71=14\"),\n  @(\"0+34=34\", \"62+35=97\"),\n  @(\"21+69=90\", \"18+68=86\"),\n  @(\"14+73=87\", \"70-35=35\"),\n  @(\"15+63=78\", \"87-52=35\"),\n  @(\"15+74=89\", \"28+20=48\"),\n  @(\"90+6=96\", \"74-22=52\"),\n  @(\"59+14=73\", \"0+83=83\"),\n  @(\"20+20=40\", \"66-10=56\"),\n  @(\"45+47=92\", \"3+41=44\"),\n  @(\"77-11=66\", \"65+2=67\"),\n  @(\"45-35=10\", \"98-66=32\"),\n  @(\"79-76=3\", \"83-72=11\"),\n  @(\"76+9=85\", \"3+0=3\"),\n  @(\"14+78=92\", \"97-12=85\"),\n  @(\"8+17=25\", \"3-3=0\"),\n  @(\"9+34=43\", \"6+2=8\"),\n  @(\"37-5=32\", \"11-5=6\"),\n  @(\"67-61=6\", \"39+9=48\"),\n  @(\"0+97=97\", \"73-28=45\"),\n  @(\"56-6=50\", \"75-4=71\"),\n  @(\"67-16=51\", \"6+57=63\"),\n  @(\"48-12=36\", \"46+3=49\"),\n  @(\"1+98=99\", \"11+83=94\"),\n  @(\"37-22=15\", \"33-31=2\"),\n  @(\"83-20=63\", \"89-23=66\"),\n  @(\"40-35=5\", \"36+41=77\"),\n  @(\"31-20=11\", \"12+50=62\"),\n  @(\"58-52=6\", \"6+14=20\"),\n  @(\"48+48=96\", \"92-79=13\"),\n  @(\"14+15=29\", \"66-19=47\"),\n  @(\"49-49=0\", \"97-37=60\"),\n  @(\"37+62=99\", \"27+20=47\"),\n  @(\"51+48=99\", \"42+31=73\"),\n  @(\"37-13=24\", \"77-24=53\"),\n  @(\"73-69=4\", \"21+23=44\"),\n  @(\"74-55=19\", \"75+13=88\"),\n  @(\"21-9=12\", \"94-43=51\"),\n  @(\"63-50=13\", \"46+44=90\"),\n  @(\"42-18=24\", \"57+42=99\"),\n  @(\"96-2=94\", \"58-58=0\"),\n  @(\"56-39=17\", \"84-0=84\"),\n  @(\"84-14=70\", \"73-0=73\"),\n  @(\"26+73=99\", \"70+29=99\"),\n  @(\"65-45=20\", \"48+31=79\"),\n  @(\"25+24=49\", \"52-46=6\"),\n  @(\"53+39=92\", \"90-31=59\"),\n  @(\"75-34=41\", \"70-39=31\"),\n  @(\"96-36=60\", \"75-52=23\"),\n  @(\"48+6=54\", \"97-43=54\"),\n  @(\"99-64=35\", \"25+0=25\"),\n  @(\"32+31=63\", \"21+16=37\"),\n  @(\"99-78=21\", \"90-68=22\"),\n  @(\"52-18=34\", \"43-11=32\"),\n  @(\"35-33=2\", \"98-23=75\"),\n  @(\"80-61=19\", \"28+28=56\"),\n  @(\"88-43=45\", \"22+11=33\"),\n  @(\"31+68=99\", \"45-1=44\"),\n  @(\"72-26=46\", \"1+62=63\"),\n  @(\"87-31=56\", \"56-38=18\"),\n  @(\"43+3=46\", \"92-90=2\"),\n  @(\"35+31=66\", \"92-40=52\"),\n  @(\"76+7=83\", \"11+68=79\"),\n  @(\"50+32=82\", \"76-3=73\"),\n  @(\"23+7=30\", \"43-27=16\"),\n  @(\"74+6=80\", \"27-26=1\"),\n  @(\"59-45=14\", \"11+42=53\"),\n  @(\"17+17=34\", \"43+7=50\"),\n  @(\"89-52=37\", \"69-20=49\"),\n  @(\"5+70=75\", \"65+24=89\"),\n  @(\"7+30=37\", \"64+25=89\"),\n  @(\"27+47=74\", \"21-19=2\"),\n  @(\"62+5=67\", \"71+0=71\"),\n  @(\"65-46=19\", \"80-8=72\"),\n  @(\"1+85=86\", \"51+26=77\"),\n  @(\"86-6=80\", \"59-59=0\"),\n  @(\"92-67=25\", \"71-19=52\"),\n  @(\"13+32=45\", \"52+33=85\"),\n  @(\"43-35=8\", \"52+19=71\"),\n  @(\"80-3=77\", \"31+30=61\"),\n  @(\"68-55=13\", \"37+23=60\"),\n  @(\"82+2=84\", \"65-12=53\"),\n  @(\"19+72=91\", \"80-45=35\"),\n  @(\"78-15=63\", \"72-71=1\"),\n  @(\"13+37=50\", \"38+51=89\"),\n  @(\"89+10=99\", \"96-25=71\"),\n  @(\"38+10=48\", \"93-56=37\"),\n  @(\"67-13=54\", \"61+18=79\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    if ($idx -ge $replacements.Count) { break }\n    $cell = $t.Cell($r, $c)\n    $cellRange = $cell.Range\n    $cellRange.MoveEnd(1, -1) | Out-Null\n    $pair = $replacements[$idx]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    if ($cellRange.Text -eq $oldText) {\n      $cellRange.Text = $newText\n    }\n    $idx += 1\n  }\n}\n"}
